$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value = 638
$ws1.Range("F3").Value = 2036
$ws1.Range("F4").Value = 47
$ws1.Range("F5").Value = 340
$ws1.Range("F6").Value = 417
$ws1.Range("F7").Value = 239
$ws1.Range("F8").Value = 13220
$ws1.Range("F10").Value = 43
$ws1.Range("F11").Value = 5354
$ws1.Range("F12").Value = 559
$ws1.Range("F13").Value = 27
$ws1.Range("F14").Value = 21
$ws1.Range("F15").Value = 40
$ws1.Range("F16").Value = 1207
$ws1.Range("F17").Value = 50
$ws1.Range("F18").Value = 142
$ws1.Range("F19").Value = 703
$ws1.Range("F20").Value = 2874
$ws1.Range("F21").Value = 7092
$ws1.Range("F22").Value = 1172
$ws1.Range("F23").Value = 3652

# --- Sheet 2: 演出 ---
$ws2 = $wb.Worksheets.Item(2)
$ws2.Range("F2").Value = 31

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value = 638
$ws4.Range("F3").Value = 2036
$ws4.Range("F4").Value = 47
$ws4.Range("F5").Value = 340
$ws4.Range("F6").Value = 31
$ws4.Range("F7").Value = 417
$ws4.Range("F8").Value = 239
$ws4.Range("F9").Value = 13220
$ws4.Range("F11").Value = 43
$ws4.Range("F12").Value = 5354
$ws4.Range("F13").Value = 559
$ws4.Range("F14").Value = 27
$ws4.Range("F15").Value = 21
$ws4.Range("F16").Value = 40
$ws4.Range("F17").Value = 1207
$ws4.Range("F18").Value = 50
$ws4.Range("F19").Value = 142
$ws4.Range("F20").Value = 703
$ws4.Range("F21").Value = 2874
$ws4.Range("F23").Value = 7092
$ws4.Range("F24").Value = 1172
$ws4.Range("F25").Value = 3652
